# Applies the changes described in the commit diff:
# 1) Update the "Expected Results" text for step 1.0 in every test case
#    (TC1..TC5) to mention ordering by number of "diarias".
# 2) Swap the Step-2 content between TC3 and TC4: TC3 becomes the
#    "atribuir/desatribuir" step (previously under TC4) and TC4 becomes the
#    "realizar a autorizacao de pagamento" step (previously under TC3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the repeated "Expected Results" text for step 1 in every TC ---
$newStep1Expected = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo numero de diarias em ordem crescente."

$ws.Range("D10").Value = $newStep1Expected
$ws.Range("D19").Value = $newStep1Expected
$ws.Range("D27").Value = $newStep1Expected
$ws.Range("D35").Value = $newStep1Expected
$ws.Range("D43").Value = $newStep1Expected

# --- 2) Swap TC3 / TC4 step-2 content ---
$tc3Step = "Chefe Clica para realizar a autorização de pagamento."
$tc3Expected = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"

$tc4Step = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$tc4Expected = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# TC3 (row 28) now gets the content that used to belong to TC4 (row 36)
$ws.Range("B28").Value = $tc4Step
$ws.Range("D28").Value = $tc4Expected

# TC4 (row 36) now gets the content that used to belong to TC3 (row 28)
$ws.Range("B36").Value = $tc3Step
$ws.Range("D36").Value = $tc3Expected
